# Update the 2D training schedule table on Sheet1.
# New data (rows 2-6, columns A-J) per updated schedule, plus an additional
# 5th trial row (row 6) that did not exist before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(1, 2, 4, 3, 9, 1, 5, 21, 5, "train_dim2_1"),
    @(2, 0, 4, 2, 8, 2, 4, 32, 5, "train_dim2_1"),
    @(3, 4, 0, 9, 1, 5, 1, 65, 5, "train_dim2_1"),
    @(4, 3, 3, 6, 6, 3, 3, 43, 5, "train_dim2_1"),
    @(5, 1, 1, 5, 3, 4, 2, 54, 5, "train_dim2_1")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}

# Selection moves to I1 in the saved file.
$ws.Range("I1").Select()
